$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 1084.2
$ws.Range("I20").Value = 1084.2
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 1084.2
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -854.2
$ws.Range("N20").ClearContents()

$ws.Range("H35").Value = 1084.2
$ws.Range("I35").Value = 1084.2
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 1084.2
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -705.2
$ws.Range("N35").ClearContents()

$ws.Range("H58").Value = 257.5
$ws.Range("I58").Value = 109
$ws.Range("J58").Value = 1000
$ws.Range("K58").Value = 327
$ws.Range("L58").Value = 3000
$ws.Range("M58").Value = -177
$ws.Range("N58").Value = -3300

$ws.Range("H64").Value = 4810930
$ws.Range("I64").Value = 6100536
$ws.Range("J64").Value = 4216.364
$ws.Range("K64").Value = 6100536
$ws.Range("L64").Value = 4216.364
$ws.Range("M64").Value = -6100288
$ws.Range("N64").Value = -4712.364

$ws.Range("H67").Value = 4810930
$ws.Range("I67").Value = 6100536
$ws.Range("J67").Value = 4216.364
$ws.Range("K67").Value = 6100536
$ws.Range("L67").Value = 4216.364
$ws.Range("M67").Value = -6099678
$ws.Range("N67").Value = -5932.364

$ws.Range("H100").Value = 4837.5
$ws.Range("I100").Value = 3950
$ws.Range("J100").Value = 5725
$ws.Range("K100").Value = 3950
$ws.Range("L100").Value = 5725
$ws.Range("M100").Value = -3409
$ws.Range("N100").Value = -6807

$ws.Range("H111").Value = 166669230
$ws.Range("I111").Value = 220
$ws.Range("J111").Value = 250003740
$ws.Range("K111").Value = 660
$ws.Range("L111").Value = 750011220
$ws.Range("M111").Value = 2407
$ws.Range("N111").Value = -750017354

$ws.Range("H113").Value = 8460.5
$ws.Range("I113").Value = 3488.5557
$ws.Range("J113").Value = 10815.632
$ws.Range("K113").Value = 3488.5557
$ws.Range("L113").Value = 10815.632
$ws.Range("M113").Value = -234.5556999999999
$ws.Range("N113").Value = -17323.632

$ws.Range("H116").Value = 3043.5715
$ws.Range("I116").Value = 2661
$ws.Range("J116").Value = 4000
$ws.Range("K116").Value = 2661
$ws.Range("L116").Value = 4000
$ws.Range("M116").Value = 781
$ws.Range("N116").Value = -10884

$ws.Range("H125").Value = 2373446.8
$ws.Range("I125").Value = 100000
$ws.Range("J125").Value = 2941808.5
$ws.Range("K125").Value = 900000
$ws.Range("L125").Value = 26476276.5
$ws.Range("M125").Value = -897540
$ws.Range("N125").Value = -26481196.5

$ws.Range("H132").Value = 38863.355
$ws.Range("I132").Value = 69798.266
$ws.Range("J132").Value = 3169.2307
$ws.Range("K132").Value = 209394.798
$ws.Range("L132").Value = 9507.6921
$ws.Range("M132").Value = -206864.798
$ws.Range("N132").Value = -14567.6921


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6884.915
$ws.Range("I32").Value = 5296.1685
$ws.Range("J32").Value = 18872.727
$ws.Range("K32").Value = 5296.1685
$ws.Range("L32").Value = 18872.727
$ws.Range("M32").Value = -5009.1685
$ws.Range("N32").Value = -19446.727

$ws.Range("H45").Value = 1319.5
$ws.Range("I45").Value = 1183.4
$ws.Range("J45").Value = 2000
$ws.Range("K45").Value = 1183.4
$ws.Range("L45").Value = 2000
$ws.Range("M45").Value = -806.4000000000001

$ws.Range("H63").Value = 3142.1428
$ws.Range("I63").Value = 2999
$ws.Range("J63").Value = 3500
$ws.Range("K63").Value = 2999
$ws.Range("L63").Value = 3500
$ws.Range("M63").Value = -2313

$ws.Range("H66").Value = 3142.1428
$ws.Range("I66").Value = 2999
$ws.Range("J66").Value = 3500
$ws.Range("K66").Value = 14995
$ws.Range("L66").Value = 17500
$ws.Range("M66").Value = -11563

$ws.Range("H110").Value = 63072.75
$ws.Range("I110").Value = 83930.336
$ws.Range("J110").Value = 500
$ws.Range("K110").Value = 83930.336
$ws.Range("L110").Value = 500
$ws.Range("M110").Value = -81885.336
$ws.Range("N110").Value = -4590


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H24").Value = 300
$ws.Range("I24").Value = 300
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 300
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -65

$ws.Range("H25").Value = 30000
$ws.Range("I25").Value = 30000
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 30000
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -29765
$ws.Range("N25").ClearContents()

$ws.Range("H29").Value = 945
$ws.Range("I29").Value = 534
$ws.Range("J29").Value = 3000
$ws.Range("K29").Value = 534
$ws.Range("L29").Value = 3000
$ws.Range("M29").Value = -245
$ws.Range("N29").Value = -3578


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1655.5714
$ws.Range("I22").Value = 1655.5714
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 1655.5714
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -1305.5714
$ws.Range("N22").ClearContents()

$ws.Range("H23").Value = 6500.5
$ws.Range("I23").Value = 6500.5
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 6500.5
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -6260.5
$ws.Range("N23").ClearContents()

$ws.Range("H27").Value = 6500.5
$ws.Range("I27").Value = 6500.5
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 6500.5
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -6308.5
$ws.Range("N27").ClearContents()

$ws.Range("H86").Value = 23813958
$ws.Range("I86").Value = 31253712
$ws.Range("J86").Value = 6740
$ws.Range("K86").Value = 31253712
$ws.Range("L86").Value = 6740
$ws.Range("M86").Value = -31252589
$ws.Range("N86").Value = -8986

$ws.Range("H89").Value = 23813958
$ws.Range("I89").Value = 31253712
$ws.Range("J89").Value = 6740
$ws.Range("K89").Value = 156268560
$ws.Range("L89").Value = 33700
$ws.Range("M89").Value = -156262944
$ws.Range("N89").Value = -44932


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 2299.0908
$ws.Range("I19").Value = 290
$ws.Range("J19").Value = 2500
$ws.Range("K19").Value = 870
$ws.Range("L19").Value = 7500
$ws.Range("M19").Value = -696
$ws.Range("N19").Value = -7848

$ws.Range("H68").Value = 4011345.2
$ws.Range("I68").Value = 13333833
$ws.Range("J68").Value = 15993.286
$ws.Range("K68").Value = 40001499
$ws.Range("L68").Value = 47979.858
$ws.Range("M68").Value = -40000688
$ws.Range("N68").Value = -49601.858

$ws.Range("H71").Value = 4011345.2
$ws.Range("I71").Value = 13333833
$ws.Range("J71").Value = 15993.286
$ws.Range("K71").Value = 120004497
$ws.Range("L71").Value = 143939.574
$ws.Range("M71").Value = -120000441
$ws.Range("N71").Value = -152051.574

$ws.Range("H122").Value = 941.8214
$ws.Range("I122").Value = 494.88
$ws.Range("J122").Value = 4666.3335
$ws.Range("K122").Value = 4453.92
$ws.Range("L122").Value = 41997.0015
$ws.Range("M122").Value = -2003.92
$ws.Range("N122").Value = -46897.0015


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 1422.3334
$ws.Range("I9").Value = 350.125
$ws.Range("J9").Value = 10000
$ws.Range("K9").Value = 350.125
$ws.Range("L9").Value = 10000
$ws.Range("M9").Value = -126.125
$ws.Range("N9").Value = -10448

$ws.Range("H22").Value = 695.5
$ws.Range("I22").Value = 646.25
$ws.Range("J22").Value = 728.3333
$ws.Range("K22").Value = 646.25
$ws.Range("L22").Value = 728.3333
$ws.Range("M22").Value = -351.25
$ws.Range("N22").Value = -1318.3333

$ws.Range("H27").Value = 695.5
$ws.Range("I27").Value = 646.25
$ws.Range("J27").Value = 728.3333
$ws.Range("K27").Value = 646.25
$ws.Range("L27").Value = 728.3333
$ws.Range("M27").Value = -539.25
$ws.Range("N27").Value = -942.3333

$ws.Range("H61").Value = 1333.3334
$ws.Range("I61").Value = 1000
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 1000
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -798
$ws.Range("N61").Value = -2404

$ws.Range("H82").Value = 1317
$ws.Range("I82").Value = 1411.6666
$ws.Range("J82").Value = 1175
$ws.Range("K82").Value = 1411.6666
$ws.Range("L82").Value = 1175
$ws.Range("M82").Value = -1050.6666
$ws.Range("N82").Value = -1897

$ws.Range("H85").Value = 1317
$ws.Range("I85").Value = 1411.6666
$ws.Range("J85").Value = 1175
$ws.Range("K85").Value = 1411.6666
$ws.Range("L85").Value = 1175
$ws.Range("M85").Value = -163.6666
$ws.Range("N85").Value = -3671

$ws.Range("H93").Value = 697.619
$ws.Range("I93").Value = 659.44446
$ws.Range("J93").Value = 926.6667
$ws.Range("K93").Value = 659.44446
$ws.Range("L93").Value = 926.6667
$ws.Range("M93").Value = 588.55554
$ws.Range("N93").Value = -3422.6667

$ws.Range("H113").Value = 1333.3334
$ws.Range("I113").Value = 1000
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1000
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 1170
$ws.Range("N113").Value = -6340


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 28750
$ws.Range("I22").Value = 5000
$ws.Range("J22").Value = 36666.668
$ws.Range("K22").Value = 5000
$ws.Range("L22").Value = 36666.668
$ws.Range("M22").Value = -4707
$ws.Range("N22").Value = -37252.668

$ws.Range("H31").Value = 21063.8
$ws.Range("I31").Value = 1833.3334
$ws.Range("J31").Value = 49909.5
$ws.Range("K31").Value = 1833.3334
$ws.Range("L31").Value = 49909.5
$ws.Range("M31").Value = -1485.3334
$ws.Range("N31").Value = -50605.5

